$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1722.125
$ws.Range("I19").Value = 798.8889
$ws.Range("K19").Value = 798.8889
$ws.Range("M19").Value = -623.8889
$ws.Range("H64").Value = 7149.875
$ws.Range("I64").Value = 6739.8
$ws.Range("K64").Value = 6739.8
$ws.Range("M64").Value = -6491.8
$ws.Range("H67").Value = 7149.875
$ws.Range("I67").Value = 6739.8
$ws.Range("K67").Value = 6739.8
$ws.Range("M67").Value = -5881.8
$ws.Range("H76").Value = 7596
$ws.Range("I76").Value = 7326.6665
$ws.Range("K76").Value = 7326.6665
$ws.Range("M76").Value = -7011.6665
$ws.Range("H79").Value = 7596
$ws.Range("I79").Value = 7326.6665
$ws.Range("K79").Value = 7326.6665
$ws.Range("M79").Value = -6234.6665
$ws.Range("H80").Value = 1000.8333
$ws.Range("J80").Value = 1189.3334
$ws.Range("L80").Value = 3568.0002
$ws.Range("N80").Value = -5564.0002
$ws.Range("H83").Value = 1000.8333
$ws.Range("J83").Value = 1189.3334
$ws.Range("L83").Value = 10704.0006
$ws.Range("N83").Value = -20688.0006
$ws.Range("H86").Value = 2814.6667
$ws.Range("I86").Value = 2966.6667
$ws.Range("J86").Value = 2662.6667
$ws.Range("K86").Value = 2966.6667
$ws.Range("L86").Value = 2662.6667
$ws.Range("M86").Value = -1843.6667
$ws.Range("N86").Value = -4908.6667
$ws.Range("H89").Value = 2814.6667
$ws.Range("I89").Value = 2966.6667
$ws.Range("J89").Value = 2662.6667
$ws.Range("K89").Value = 14833.3335
$ws.Range("L89").Value = 13313.3335
$ws.Range("M89").Value = -9217.333500000001
$ws.Range("N89").Value = -24545.3335
$ws.Range("H92").Value = 769.931
$ws.Range("I92").Value = 782.2692
$ws.Range("J92").Value = 663
$ws.Range("K92").Value = 782.2692
$ws.Range("L92").Value = 663
$ws.Range("M92").Value = 465.7308
$ws.Range("N92").Value = -3159
$ws.Range("H96").Value = 792
$ws.Range("I96").Value = 683.1667
$ws.Range("J96").Value = 900.8333
$ws.Range("K96").Value = 2049.5001
$ws.Range("L96").Value = 2702.4999
$ws.Range("M96").Value = -676.5001000000002
$ws.Range("N96").Value = -5448.4999
$ws.Range("H106").Value = 4982
$ws.Range("I106").Value = 4980.8335
$ws.Range("K106").Value = 4980.8335
$ws.Range("M106").Value = -4349.8335
$ws.Range("H111").Value = 523.3333
$ws.Range("J111").Value = 470.5
$ws.Range("L111").Value = 1411.5
$ws.Range("N111").Value = -7545.5
$ws.Range("H132").Value = 4710.1904
$ws.Range("I132").Value = 5273.0557
$ws.Range("K132").Value = 15819.1671
$ws.Range("M132").Value = -13289.1671
$ws.Range("H137").Value = 60360.453
$ws.Range("I137").Value = 107227.65
$ws.Range("K137").Value = 321682.95
$ws.Range("M137").Value = -319132.95

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3965.04
$ws.Range("I32").Value = 1908.9193
$ws.Range("J32").Value = 13771.154
$ws.Range("K32").Value = 1908.9193
$ws.Range("L32").Value = 13771.154
$ws.Range("M32").Value = -1621.9193
$ws.Range("N32").Value = -14345.154
$ws.Range("H74").Value = 56880.35
$ws.Range("I74").Value = 4052.2222
$ws.Range("J74").Value = 260646
$ws.Range("K74").Value = 4052.2222
$ws.Range("L74").Value = 260646
$ws.Range("M74").Value = -3178.2222
$ws.Range("N74").Value = -262394
$ws.Range("H77").Value = 56880.35
$ws.Range("I77").Value = 4052.2222
$ws.Range("J77").Value = 260646
$ws.Range("K77").Value = 20261.111
$ws.Range("L77").Value = 1303230
$ws.Range("M77").Value = -15893.111
$ws.Range("N77").Value = -1311966
$ws.Range("H88").Value = 1520.6154
$ws.Range("I88").Value = 681.25
$ws.Range("J88").Value = 1893.6666
$ws.Range("K88").Value = 681.25
$ws.Range("L88").Value = 1893.6666
$ws.Range("M88").Value = -275.25
$ws.Range("N88").Value = -2705.6666
$ws.Range("H91").Value = 1520.6154
$ws.Range("I91").Value = 681.25
$ws.Range("J91").Value = 1893.6666
$ws.Range("K91").Value = 681.25
$ws.Range("L91").Value = 1893.6666
$ws.Range("M91").Value = 722.75
$ws.Range("N91").Value = -4701.6666
$ws.Range("H102").Value = 4131.4736
$ws.Range("I102").Value = 3245.2666
$ws.Range("K102").Value = 3245.2666
$ws.Range("M102").Value = -1623.2666
$ws.Range("H122").Value = 11114016
$ws.Range("I122").Value = 13891145
$ws.Range("K122").Value = 41673435
$ws.Range("M122").Value = -41670985
$ws.Range("H132").Value = 1466.7872
$ws.Range("I132").Value = 878.44446
$ws.Range("J132").Value = 3392.2727
$ws.Range("K132").Value = 2635.33338
$ws.Range("L132").Value = 10176.8181
$ws.Range("M132").Value = -105.33338
$ws.Range("N132").Value = -15236.8181

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 5648.2104
$ws.Range("I94").Value = 1139.909
$ws.Range("K94").Value = 1139.909
$ws.Range("M94").Value = -688.9090000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 49431.316
$ws.Range("I132").Value = 2397.375
$ws.Range("K132").Value = 7192.125
$ws.Range("M132").Value = -4662.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 57.25
$ws.Range("I10").Value = 73
$ws.Range("K10").Value = 219
$ws.Range("M10").Value = -80
$ws.Range("H105").Value = 4836.75
$ws.Range("J105").Value = 4836.75
$ws.Range("L105").Value = 14510.25
$ws.Range("N105").Value = -19752.25
$ws.Range("H134").Value = 3900
$ws.Range("I134").Value = 3900
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 11700
$ws.Range("L134").Value = 0
$ws.Range("N134").Value = -6630
$ws.Range("H140").Value = 2468.5881
$ws.Range("I140").Value = 2264.4666
$ws.Range("K140").Value = 6793.399800000001
$ws.Range("M140").Value = -1613.399800000001
$ws.Range("H141").Value = 2357.8333
$ws.Range("I141").Value = 1965.6666
$ws.Range("K141").Value = 5896.9998
$ws.Range("M141").Value = -716.9997999999996

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 9115
$ws.Range("I80").Value = 3300
$ws.Range("K80").Value = 3300
$ws.Range("M80").Value = -2302
$ws.Range("H83").Value = 9115
$ws.Range("I83").Value = 3300
$ws.Range("K83").Value = 16500
$ws.Range("M83").Value = -11508
$ws.Range("H102").Value = 37510.45
$ws.Range("I102").Value = 2390.9
$ws.Range("J102").Value = 115553.89
$ws.Range("K102").Value = 2390.9
$ws.Range("L102").Value = 115553.89
$ws.Range("M102").Value = -768.9000000000001
$ws.Range("N102").Value = -118797.89
$ws.Range("H126").Value = 3563.7646
$ws.Range("I126").Value = 3327.8572
$ws.Range("K126").Value = 9983.571599999999
$ws.Range("M126").Value = -7513.571599999999
$ws.Range("H132").Value = 3105.1555
$ws.Range("I132").Value = 2456.6667
$ws.Range("K132").Value = 7370.000100000001
$ws.Range("M132").Value = -4840.000100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 58005.5
$ws.Range("I22").Value = 84258.336
$ws.Range("J22").Value = 5499.8335
$ws.Range("K22").Value = 84258.336
$ws.Range("L22").Value = 5499.8335
$ws.Range("M22").Value = -83963.336
$ws.Range("N22").Value = -6089.8335
$ws.Range("H27").Value = 58005.5
$ws.Range("I27").Value = 84258.336
$ws.Range("J27").Value = 5499.8335
$ws.Range("K27").Value = 84258.336
$ws.Range("L27").Value = 5499.8335
$ws.Range("M27").Value = -84151.336
$ws.Range("N27").Value = -5713.8335
$ws.Range("H41").Value = 27868.5
$ws.Range("J41").Value = 32324.666
$ws.Range("L41").Value = 32324.666
$ws.Range("N41").Value = -33200.666
$ws.Range("H46").Value = 4914.5293
$ws.Range("I46").Value = 3779.6
$ws.Range("J46").Value = 6535.857
$ws.Range("K46").Value = 3779.6
$ws.Range("L46").Value = 6535.857
$ws.Range("M46").Value = -3591.6
$ws.Range("N46").Value = -6911.857
$ws.Range("H82").Value = 1167.2
$ws.Range("I82").Value = 842.375
$ws.Range("J82").Value = 1538.4286
$ws.Range("K82").Value = 842.375
$ws.Range("L82").Value = 1538.4286
$ws.Range("M82").Value = -481.375
$ws.Range("N82").Value = -2260.4286
$ws.Range("H85").Value = 1167.2
$ws.Range("I85").Value = 842.375
$ws.Range("J85").Value = 1538.4286
$ws.Range("K85").Value = 842.375
$ws.Range("L85").Value = 1538.4286
$ws.Range("M85").Value = 405.625
$ws.Range("N85").Value = -4034.4286
$ws.Range("H136").Value = 54890.4
$ws.Range("I136").Value = 80447
$ws.Range("J136").Value = 7428.143
$ws.Range("K136").Value = 241341
$ws.Range("L136").Value = 22284.429
$ws.Range("M136").Value = -238791
$ws.Range("N136").Value = -27384.429

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1433.6666
$ws.Range("I81").Value = 1525
$ws.Range("J81").Value = 1251
$ws.Range("K81").Value = 3050
$ws.Range("L81").Value = 2502
$ws.Range("M81").Value = -1989
$ws.Range("N81").Value = -4624
$ws.Range("H84").Value = 1433.6666
$ws.Range("I84").Value = 1525
$ws.Range("J84").Value = 1251
$ws.Range("K84").Value = 15250
$ws.Range("L84").Value = 12510
$ws.Range("M84").Value = -9946
$ws.Range("N84").Value = -23118
$ws.Range("H113").Value = 830.04346
$ws.Range("I113").Value = 611.9286
$ws.Range("K113").Value = 1835.7858
$ws.Range("M113").Value = 334.2142000000001
$ws.Range("H126").Value = 3680.1333
$ws.Range("I126").Value = 3764.4285
$ws.Range("K126").Value = 11293.2855
$ws.Range("M126").Value = -8823.2855
$ws.Range("H136").Value = 3238.0344
$ws.Range("I136").Value = 2652.7368
$ws.Range("J136").Value = 4350.1
$ws.Range("K136").Value = 7958.2104
$ws.Range("L136").Value = 13050.3
$ws.Range("M136").Value = -5408.2104
$ws.Range("N136").Value = -18150.3
